# "Alter bei AN-Anteil PV hinzugefuegt und notwendige Aenderungen vorgenommen."
#
# Insert a new question row ("juenger als 23 oder vor 1940 geboren?") right
# above "wohnhaft Sachsen?" on Tabelle1, and blank out the answers for the
# three yes/no/number questions that now need to be re-answered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row before the current row 41 ("wohnhaft Sachsen?").
# This pushes row 41 onward down by one and carries formatting/validation
# with it, matching rows 39/40 staying put and rows 41.. shifting to 42..
$ws.Rows.Item(41).Insert()

# New question text in column A of the freshly inserted row.
$ws.Range("A41").Value = "juenger als 23 oder vor 1940 geboren?"

# Clear the previously filled-in answers that need to be reconsidered:
#  - B39 "ermaessigter Krankenversicherungsbeitrag?" (was "nein")
#  - B40 "Anzahl Kinder" (was 1)
#  - B42 "wohnhaft Sachsen?" (was "ja") - this is the row that moved from 41 to 42
$ws.Range("B39").ClearContents()
$ws.Range("B40").ClearContents()
$ws.Range("B41").ClearContents()
$ws.Range("B42").ClearContents()

# Update the view: scroll down to the area being edited and select B40.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("B40").Select()
